$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Q1 for "Business Reporter" row: job titles question rewording
$ws.Range("A2").Value = "What is each business person’s job title?"

# Q3 for "Restaurant Critic" row: clients order question rewording (keeps trailing blank line)
$ws.Range("C3").Value = "What does each client order?`n"

# D1 for "Wedding Planner" row: couple dating duration wording
$ws.Range("A5").Value = "For how long has the couple been dating?"

# D2 for "Wedding Planner" row: diamonds/ring color question - add comma, drop trailing newline
$ws.Range("B5").Value = "How many diamonds are on the ring, and what is the diamond color?"

# Q2 for "Grocery Store Customer Experience Manager" row: pick out -> pick out to buy, drop trailing newline
$ws.Range("B7").Value = "What items do the clients pick out to buy?"

# Q3 for "Grocery Store Customer Experience Manager" row: add comma
$ws.Range("C7").Value = "How many checkout lanes are open, and which one do the clients step into?`n"

# Q4 for "Grocery Store Customer Experience Manager" row: add comma
$ws.Range("D7").Value = "How much are the groceries, and what method of payment do the clients use?`n"

# Q2 for "Dean of Academic Studies" row: add comma
$ws.Range("B8").Value = "What class are the students in, and what is the day’s lecture about?`n"

# Q4 for "Dean of Academic Studies" row: "in class" -> "for the class", add comma
$ws.Range("D8").Value = "What is the next assessment/assignment for the class, and when is it scheduled/due?"

# D1 for "Couples Therapist" row: add comma, add trailing newline
$ws.Range("C9").Value = "Does the person who is being broken up with want to break up, and what’s the reason stated by the person being broken up with that he/she does or does not want to break up?`n"

# D2 for "Couples Therapist" row: reworded question about items back
$ws.Range("D9").Value = "Who wants what items back as a result of the breakup?"
